# Lecture 26.pptx edit:
#   Slide 17, "Content Placeholder 2" shape contains a small Java
#   switch-expression code sample.  The "Wednesday"/"Saturday" weekday
#   cases (and their "yield" lines) are removed, merging directly into
#   the "Sunday" case a few lines down.  The remaining "Sunday" line
#   keeps its original four leading spaces, just now carried by a
#   dedicated leading run ("    ") followed by the case text itself.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(17)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Paragraphs, joined with CR (PowerPoint's paragraph separator).
$lines = @(
    'return  switch (day) {',
    '    case "Monday":',
    '        yield  "Weekday";',
    '    case "Tuesday":',
    '        yield "Weekday";',
    '    case "Sunday":',
    '        yield "Weekend";',
    '    default:',
    '        yield "Unknown";',
    '};'
)
$tr.Text = [string]::Join("`r", $lines)

# Split the "    case "Sunday":" paragraph's single run into two runs —
# a leading 4-space run and the "case "Sunday":" text run — by nudging
# a (no-op) character-level format on just the leading spaces.
$sundayPara = $tr.Paragraphs(6, 1)
$leadingSpaces = $sundayPara.Characters(1, 4)
$leadingSpaces.Font.Size = $leadingSpaces.Font.Size
